# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) values in row 2 of the
# zh-cn and de-de worksheets to reflect the new report-generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 12:06:55"
$wsZhCn.Range("H2").Value = "2016-03-24 12:07:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 12:06:59"
$wsDeDe.Range("H2").Value = "2016-03-24 12:07:38"
